{"js": "// Replace the arithmetic-expression text in every cell of the first table\n// with its updated value, in row-major order, while preserving the\n// existing paragraph / run formatting (font, size, justification, etc.).\nconst pairs = [[\"9+40=49\", \"34-21=13\"], [\"0+2=2\", \"41-31=10\"], [\"27+59=86\", \"51-4=47\"], [\"14+34=48\", \"77+13=90\"], [\"6+73=79\", \"78-75=3\"], [\"67-16=51\", \"48+6=54\"], [\"65-63=2\", \"34+13=47\"], [\"93-5=88\", \"44+25=69\"], [\"7+49=56\", \"13+76=89\"], [\"0+40=40\", \"43+7=50\"], [\"85-19=66\", \"34+37=71\"], [\"52+39=91\", \"79-14=65\"], [\"31+30=61\", \"71-2=69\"], [\"25+39=64\", \"50-18=32\"], [\"12+17=29\", \"52-27=25\"], [\"24+74=98\", \"35+24=59\"], [\"25-11=14\", \"66+6=72\"], [\"79-22=57\", \"59-15=44\"], [\"60-3=57\", \"35+13=48\"], [\"9-2=7\", \"81-70=11\"], [\"31+0=31\", \"36+11=47\"], [\"6+37=43\", \"65-49=16\"], [\"77-1=76\", \"28+18=46\"], [\"88-57=31\", \"51+11=62\"], [\"66+33=99\", \"41+27=68\"], [\"96-46=50\", \"6+77=83\"], [\"99-54=45\", \"51-0=51\"], [\"16+78=94\", \"47-40=7\"], [\"18+22=40\", \"0+59=59\"], [\"79+20=99\", \"36+59=95\"], [\"52-6=46\", \"6+56=62\"], [\"12+19=31\", \"90+4=94\"], [\"72-6=66\", \"77-15=62\"], [\"24-14=10\", \"82-15=67\"], [\"41-28=13\", \"87-18=69\"], [\"0+39=39\", \"74-25=49\"], [\"10+50=60\", \"43+40=83\"], [\"84-29=55\", \"20+46=66\"], [\"9+17=26\", \"86-7=79\"], [\"5+9=14\", \"78+19=97\"], [\"45+44=89\", \"14+41=55\"], [\"65-63=2\", \"14+1=15\"], [\"7+28=35\", \"92-38=54\"], [\"6+46=52\", \"44+39=83\"], [\"53-6=47\", \"53+17=70\"], [\"51-35=16\", \"35+8=43\"], [\"85-76=9\", \"82-54=28\"], [\"73-42=31\", \"35+41=76\"], [\"62-9=53\", \"44-6=38\"], [\"51+22=73\", \"87-1=86\"], [\"21+58=79\", \"66+33=99\"], [\"69-4=65\", \"93+4=97\"], [\"93-61=32\", \"90-28=62\"], [\"99-10=89\", \"58+20=78\"], [\"8+83=91\", \"35-24=11\"], [\"59-45=14\", \"50-12=38\"], [\"22-8=14\", \"54+42=96\"], [\"8+84=92\", \"64-46=18\"], [\"83-20=63\", \"50+41=91\"], [\"43-18=25\", \"47-29=18\"], [\"90-70=20\", \"48+7=55\"], [\"93-10=83\", \"8+56=64\"], [\"21+70=91\", \"64-23=41\"], [\"50-22=28\", \"89-71=18\"], [\"1+95=96\", \"68-31=37\"], [\"16+58=74\", \"36-11=25\"], [\"80-1=79\", \"99-74=25\"], [\"32+30=62\", \"63-50=13\"], [\"32+60=92\", \"1+35=36\"], [\"3+87=90\", \"7+0=7\"], [\"91-42=49\", \"79-38=41\"], [\"21+55=76\", \"78+9=87\"], [\"95-73=22\", \"91-79=12\"], [\"24+28=52\", \"92-9=83\"], [\"42+18=60\", \"20+44=64\"], [\"65-47=18\", \"6+26=32\"], [\"89-46=43\", \"73-53=20\"], [\"84-46=38\", \"60+32=92\"], [\"27+33=60\", \"54+31=85\"], [\"74+4=78\", \"98-26=72\"], [\"71-18=53\", \"49+21=70\"], [\"8+9=17\", \"40-24=16\"], [\"6+29=35\", \"9+11=20\"], [\"92-63=29\", \"84+14=98\"], [\"55-7=48\", \"79+5=84\"], [\"47-39=8\", \"11+23=34\"], [\"49+46=95\", \"97-28=69\"], [\"71-19=52\", \"13+44=57\"], [\"95-50=45\", \"96-61=35\"], [\"22+20=42\", \"23+50=73\"], [\"9+70=79\", \"47+32=79\"], [\"8+24=32\", \"0+45=45\"], [\"22-4=18\", \"49+0=49\"], [\"5+36=41\", \"45+54=99\"], [\"37+47=84\", \"51-26=25\"], [\"71-24=47\", \"44-36=8\"], [\"75-40=35\", \"7+62=69\"], [\"17+50=67\", \"23+42=65\"], [\"34+23=57\", \"22+50=72\"], [\"84-10=74\", \"44-10=34\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values.length ? table.values[0].length : 0;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount && idx < pairs.length; r++) {\n  for (let c = 0; c < colCount && idx < pairs.length; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    // eslint-disable-next-line no-await-in-loop\n    await context.sync();\n\n    const [oldVal, newVal] = pairs[idx];\n    if (para.text === oldVal) {\n      para.getRange().insertText(newVal, Word.InsertLocation.replace);\n    }\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic-expression text in every cell of the first table\n# with its updated value, in row-major order, while preserving the\n# existing paragraph / run formatting (font, size, justification, etc.).\n$pairs = @(\n    @(\"9+40=49\", \"34-21=13\"),\n    @(\"0+2=2\", \"41-31=10\"),\n    @(\"27+59=86\", \"51-4=47\"),\n    @(\"14+34=48\", \"77+13=90\"),\n    @(\"6+73=79\", \"78-75=3\"),\n    @(\"67-16=51\", \"48+6=54\"),\n    @(\"65-63=2\", \"34+13=47\"),\n    @(\"93-5=88\", \"44+25=69\"),\n    @(\"7+49=56\", \"13+76=89\"),\n    @(\"0+40=40\", \"43+7=50\"),\n    @(\"85-19=66\", \"34+37=71\"),\n    @(\"52+39=91\", \"79-14=65\"),\n    @(\"31+30=61\", \"71-2=69\"),\n    @(\"25+39=64\", \"50-18=32\"),\n    @(\"12+17=29\", \"52-27=25\"),\n    @(\"24+74=98\", \"35+24=59\"),\n    @(\"25-11=14\", \"66+6=72\"),\n    @(\"79-22=57\", \"59-15=44\"),\n    @(\"60-3=57\", \"35+13=48\"),\n    @(\"9-2=7\", \"81-70=11\"),\n    @(\"31+0=31\", \"36+11=47\"),\n    @(\"6+37=43\", \"65-49=16\"),\n    @(\"77-1=76\", \"28+18=46\"),\n    @(\"88-57=31\", \"51+11=62\"),\n    @(\"66+33=99\", \"41+27=68\"),\n    @(\"96-46=50\", \"6+77=83\"),\n    @(\"99-54=45\", \"51-0=51\"),\n    @(\"16+78=94\", \"47-40=7\"),\n    @(\"18+22=40\", \"0+59=59\"),\n    @(\"79+20=99\", \"36+59=95\"),\n    @(\"52-6=46\", \"6+56=62\"),\n    @(\"12+19=31\", \"90+4=94\"),\n    @(\"72-6=66\", \"77-15=62\"),\n    @(\"24-14=10\", \"82-15=67\"),\n    @(\"41-28=13\", \"87-18=69\"),\n    @(\"0+39=39\", \"74-25=49\"),\n    @(\"10+50=60\", \"43+40=83\"),\n    @(\"84-29=55\", \"20+46=66\"),\n    @(\"9+17=26\", \"86-7=79\"),\n    @(\"5+9=14\", \"78+19=97\"),\n    @(\"45+44=89\", \"14+41=55\"),\n    @(\"65-63=2\", \"14+1=15\"),\n    @(\"7+28=35\", \"92-38=54\"),\n    @(\"6+46=52\", \"44+39=83\"),\n    @(\"53-6=47\", \"53+17=70\"),\n    @(\"51-35=16\", \"35+8=43\"),\n    @(\"85-76=9\", \"82-54=28\"),\n    @(\"73-42=31\", \"35+41=76\"),\n    @(\"62-9=53\", \"44-6=38\"),\n    @(\"51+22=73\", \"87-1=86\"),\n    @(\"21+58=79\", \"66+33=99\"),\n    @(\"69-4=65\", \"93+4=97\"),\n    @(\"93-61=32\", \"90-28=62\"),\n    @(\"99-10=89\", \"58+20=78\"),\n    @(\"8+83=91\", \"35-24=11\"),\n    @(\"59-45=14\", \"50-12=38\"),\n    @(\"22-8=14\", \"54+42=96\"),\n    @(\"8+84=92\", \"64-46=18\"),\n    @(\"83-20=63\", \"50+41=91\"),\n    @(\"43-18=25\", \"47-29=18\"),\n    @(\"90-70=20\", \"48+7=55\"),\n    @(\"93-10=83\", \"8+56=64\"),\n    @(\"21+70=91\", \"64-23=41\"),\n    @(\"50-22=28\", \"89-71=18\"),\n    @(\"1+95=96\", \"68-31=37\"),\n    @(\"16+58=74\", \"36-11=25\"),\n    @(\"80-1=79\", \"99-74=25\"),\n    @(\"32+30=62\", \"63-50=13\"),\n    @(\"32+60=92\", \"1+35=36\"),\n    @(\"3+87=90\", \"7+0=7\"),\n    @(\"91-42=49\", \"79-38=41\"),\n    @(\"21+55=76\", \"78+9=87\"),\n    @(\"95-73=22\", \"91-79=12\"),\n    @(\"24+28=52\", \"92-9=83\"),\n    @(\"42+18=60\", \"20+44=64\"),\n    @(\"65-47=18\", \"6+26=32\"),\n    @(\"89-46=43\", \"73-53=20\"),\n    @(\"84-46=38\", \"60+32=92\"),\n    @(\"27+33=60\", \"54+31=85\"),\n    @(\"74+4=78\", \"98-26=72\"),\n    @(\"71-18=53\", \"49+21=70\"),\n    @(\"8+9=17\", \"40-24=16\"),\n    @(\"6+29=35\", \"9+11=20\"),\n    @(\"92-63=29\", \"84+14=98\"),\n    @(\"55-7=48\", \"79+5=84\"),\n    @(\"47-39=8\", \"11+23=34\"),\n    @(\"49+46=95\", \"97-28=69\"),\n    @(\"71-19=52\", \"13+44=57\"),\n    @(\"95-50=45\", \"96-61=35\"),\n    @(\"22+20=42\", \"23+50=73\"),\n    @(\"9+70=79\", \"47+32=79\"),\n    @(\"8+24=32\", \"0+45=45\"),\n    @(\"22-4=18\", \"49+0=49\"),\n    @(\"5+36=41\", \"45+54=99\"),\n    @(\"37+47=84\", \"51-26=25\"),\n    @(\"71-24=47\", \"44-36=8\"),\n    @(\"75-40=35\", \"7+62=69\"),\n    @(\"17+50=67\", \"23+42=65\"),\n    @(\"34+23=57\", \"22+50=72\"),\n    @(\"84-10=74\", \"44-10=34\")\n)\n\n$wdCharacter = 1\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        if ($idx -ge $pairs.Count) { break }\n\n        $oldVal = $pairs[$idx][0]\n        $newVal = $pairs[$idx][1]\n\n        $cell = $tbl.Cell($r, $c)\n        $rng = $cell.Range\n        [void]$rng.MoveEnd($wdCharacter, -1)  # trim trailing end-of-cell mark\n\n        if ($rng.Text -eq $oldVal) {\n            $rng.Text = $newVal\n        }\n\n        $idx++\n    }\n}\n"}
